$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$data = @{
    2  = @{ D = 44497; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 }
    3  = @{ D = 44452; J = 120; K = 2300; L = 2300; M = 2300; P = 2300 }
    4  = @{ D = 44203; J = 30;  K = 2000; L = 2000; M = 2000; P = 2000 }
    5  = @{ D = 44447; J = 75;  K = 2200; L = 2200; M = 2200; P = 2200 }
    6  = @{ D = 44483; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 }
    7  = @{ D = 44474; J = 20;  K = 1600; L = 1600; M = 1600; P = 1600 }
    8  = @{ D = 44496; J = 40;  K = 2200; L = 2200; M = 2200; P = 2200 }
    9  = @{ D = 44453; J = 20;  K = 2300; L = 2300; M = 2300; P = 2300 }
    10 = @{ D = 44476; J = 30;  K = 2200; L = 2200; M = 2200; P = 2200 }
    11 = @{ D = 44487; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 }
    12 = @{ D = 44473; J = 140; K = 1600; L = 1600; M = 1600; P = 1600 }
    13 = @{ D = 44484; J = 40;  K = 2200; L = 2200; M = 2200; P = 2200 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
